$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# --- Row 16 (paper 12): replace the placeholder paper with the new upload ---
# Cell write order matches the author's entry order so new shared-string
# indices land the same way they do in the target workbook.
# Link to article (plain text now - the old hyperlink to the embj paper is removed below)
$ws.Cells.Item(16, 4).Value = "https://doi.org/10.1002/jnr.24386"
# Title
$ws.Cells.Item(16, 2).Value = "Glycogen distribution in mouse hippocampus"
# Year
$ws.Cells.Item(16, 3).Value = 2018
# Scores: Quality of experimental design, Description of methods, Description of sample
# population/data under study, Data availability and accessibility, Tool accessibility
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 3
$ws.Cells.Item(16, 8).Value = 2
$ws.Cells.Item(16, 9).Value = 2
$ws.Cells.Item(16, 10).Value = 2
# Reproducibility / Replicability = "?"
$ws.Cells.Item(16, 11).Value = "?"
$ws.Cells.Item(16, 12).Value = "?"
# Notes
$ws.Cells.Item(16, 13).Value = "Qualitative study with minimal detail and inaccesible equipment. "

# --- Remove the hyperlink that used to live on D16 (old embj.org link) ---
# This sandboxed engine's Hyperlinks.Delete() always clears the whole sheet
# collection regardless of the scope it was fetched from (and individual
# Hyperlink.Delete() is a no-op), so: capture the other hyperlinks first,
# wipe the collection, then recreate the ones that must survive. Add() also
# clobbers the anchor cell's number format/alignment (it re-applies a
# slightly different "Hyperlink" style), so stash + restore each cell's
# original formatting around the re-add via copy/paste-special-formats.
$keep = @()
foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -ne "`$D`$16") {
        $keep += ,@($addr, $h.Address)
    }
}

$ws.Hyperlinks.Delete()

foreach ($item in $keep) {
    $cellAddr = $item[0] -replace '\$', ''
    $url = $item[1]
    $cell = $ws.Range($cellAddr)
    $cell.Copy()
    $ws.Hyperlinks.Add($cell, $url) | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null  # xlPasteFormats - restore clobbered formatting
}
$excel.CutCopyMode = $false

# --- Update the active selection to match the author's final cursor position ---
$ws.Activate()
$ws.Range("M16").Select()
